# "Fixed MongoDB steps. Setup MongoDB only for Broker."
#
# The MongoDB install/setup walk-through (Download and install MongoDB ->
# Create directory -> Run MongoDB -> Verify) used to sit right before the
# "Run broker-mongodb.js" slide (position 9). The author moved the
# "Run broker-mongodb.js" slide up to immediately follow the "Design &
# Implementation" divider (position 5), pushing the MongoDB setup slides
# down to positions 6-9 so the MongoDB install steps now read as setup
# specifically for the broker.

$p = $ppt.ActivePresentation

# Slide 9 ("Run broker-mongodb.js") moves to slide 5; slides 5-8 (the
# MongoDB download/create-dir/run/verify steps) shift down to 6-8-... in
# order to make room, exactly like dragging slide 9 up to slot 5 in the
# slide sorter.
$p.Slides.Item(9).MoveTo(5)
